# Update the division problems in the worksheet table to the regenerated
# set of values (commit: "Update master to output generated at 4250d90").
#
# The document has a single 5-column table; every 4th row (1, 5, 9, 13, 17)
# holds one row of "NN÷N=" problems, the rows between are blank spacer rows.
# Cells are addressed directly by (row, column) so that the replacement is
# unambiguous even though a couple of the new values coincide with old
# values used elsewhere in the table (e.g. "28÷9=" is both replaced and
# used as a replacement value).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1
$t.Cell(1, 1).Range.Text = "14÷7="
$t.Cell(1, 2).Range.Text = "91÷7="
$t.Cell(1, 3).Range.Text = "37÷4="
$t.Cell(1, 4).Range.Text = "97÷5="
$t.Cell(1, 5).Range.Text = "27÷8="

# Row 5
$t.Cell(5, 1).Range.Text = "28÷9="
$t.Cell(5, 2).Range.Text = "87÷8="
$t.Cell(5, 3).Range.Text = "84÷8="
$t.Cell(5, 4).Range.Text = "60÷2="
$t.Cell(5, 5).Range.Text = "77÷2="

# Row 9
$t.Cell(9, 1).Range.Text = "84÷4="
$t.Cell(9, 2).Range.Text = "53÷8="
$t.Cell(9, 3).Range.Text = "88÷6="
$t.Cell(9, 4).Range.Text = "33÷6="
$t.Cell(9, 5).Range.Text = "13÷8="

# Row 13
$t.Cell(13, 1).Range.Text = "82÷9="
$t.Cell(13, 2).Range.Text = "87÷6="
$t.Cell(13, 3).Range.Text = "69÷8="
$t.Cell(13, 4).Range.Text = "36÷5="
$t.Cell(13, 5).Range.Text = "50÷9="

# Row 17
$t.Cell(17, 1).Range.Text = "36÷5="
$t.Cell(17, 2).Range.Text = "83÷4="
$t.Cell(17, 3).Range.Text = "20÷5="
$t.Cell(17, 4).Range.Text = "39÷2="
$t.Cell(17, 5).Range.Text = "48÷7="
